# "Added win, lose and reset functionality"
# The Lab1 task list gained two rows describing that work earlier on
# (rows 11 "Add Lose and Win Screens" and 12 "Add Reset button"); this
# edit fills in their "Actual Time to Complete" values now that the
# tasks are done, and leaves the selection where the user's cursor
# ended up afterwards (C13, just past the bottom of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = 15
$ws.Range("C12").Value = 20

[void]$ws.Range("C13").Select()
